# The commit swaps the embedded theme XML parts ppt/theme/theme1.xml
# (the slide master's theme, originally "Integral" / "Red Violet") and
# ppt/theme/theme2.xml (the notes master's theme, originally
# "Office Theme" / "Office") with each other.
#
# The two themes are identical except for their 12 scheme colours (the
# font scheme and format scheme are byte-for-byte identical), so the
# swap is reproduced here by exchanging the 12
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink colours between the slide
# master's theme and the notes master's theme via the PowerPoint theme
# colour-scheme object model.
#
# Note: ThemeColorScheme.Item(n).RGB uses the classic VBA COLORREF
# (0x00BBGGRR) byte order, not 0x00RRGGBB.

$p = $ppt.ActivePresentation

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
$notesScheme  = $p.NotesMaster.Theme.ThemeColorScheme

# Original slide-master theme colours ("Integral" / "Red Violet"),
# index 1..12 = dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$integralColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    5326149,      # dk2      454551
    14473688,     # lt2      D8D9DC
    9514467,      # accent1  E32D91
    13381832,     # accent2  C830CC
    14460494,     # accent3  4EA6DC
    15168839,     # accent4  4775E7
    14774665,     # accent5  8971E1
    7555029,      # accent6  D54773
    2465643,      # hlink    6B9F25
    9211020       # folHlink 8C8C8C
)

# Original notes-master theme colours ("Office Theme" / "Office"),
# same index order.
$officeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

# Put the "Office Theme" colours on the slide master (was "Integral").
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}

# Put the "Integral" colours on the notes master (was "Office Theme").
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
